$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1532
$ws.Range("F3").Value = 877
$ws.Range("F4").Value = 458
$ws.Range("F5").Value = 902
$ws.Range("F6").Value = 512
$ws.Range("F7").Value = 7712
$ws.Range("F8").Value = 140
$ws.Range("F10").Value = 1939
$ws.Range("F11").Value = 5568
$ws.Range("F12").Value = 571
$ws.Range("F13").Value = 309
$ws.Range("F14").Value = 7692
$ws.Range("F15").Value = 9091
$ws.Range("F16").Value = 1151
$ws.Range("F17").Value = 908
$ws.Range("F18").Value = 4476
$ws.Range("F19").Value = 676
$ws.Range("F20").Value = 242
$ws.Range("F23").Value = 163
$ws.Range("F24").Value = 1200
$ws.Range("F25").Value = 119
$ws.Range("F26").Value = 1677
$ws.Range("F27").Value = 728
$ws.Range("F28").Value = 942
$ws.Range("F29").Value = 9
$ws.Range("F30").Value = 1884
$ws.Range("F31").Value = 340
$ws.Range("F32").Value = 2315
$ws.Range("F34").Value = 115
$ws.Range("F35").Value = 1473
$ws.Range("F36").Value = 69
$ws.Range("F38").Value = 799
$ws.Range("F39").Value = 513
$ws.Range("F40").Value = 2984
$ws.Range("F41").Value = 4118
$ws.Range("F42").Value = 193
$ws.Range("F43").Value = 45
$ws.Range("F44").Value = 426
$ws.Range("F45").Value = 514
$ws.Range("F46").Value = 15
$ws.Range("F48").Value = 177
$ws.Range("F49").Value = 4091

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 23
$ws.Range("F9").Value = 23
$ws.Range("F27").Value = 36

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5261

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1532
$ws.Range("F4").Value = 877
$ws.Range("F5").Value = 458
$ws.Range("F6").Value = 902
$ws.Range("F7").Value = 512
$ws.Range("F9").Value = 23
$ws.Range("F11").Value = 5568
$ws.Range("F12").Value = 571
$ws.Range("F13").Value = 7692
$ws.Range("F15").Value = 1151
$ws.Range("F16").Value = 908
$ws.Range("F17").Value = 676
$ws.Range("F18").Value = 242
$ws.Range("F22").Value = 163
$ws.Range("F23").Value = 1200
$ws.Range("F24").Value = 119
$ws.Range("F25").Value = 1677
$ws.Range("F26").Value = 728
$ws.Range("F27").Value = 942
$ws.Range("F28").Value = 9
$ws.Range("F29").Value = 1884
$ws.Range("F30").Value = 340
$ws.Range("F31").Value = 2315
$ws.Range("F33").Value = 69
$ws.Range("F39").Value = 513
$ws.Range("F40").Value = 4118
$ws.Range("F41").Value = 36
$ws.Range("F42").Value = 193
$ws.Range("F43").Value = 45
$ws.Range("F44").Value = 426
$ws.Range("F45").Value = 514
$ws.Range("F46").Value = 15
$ws.Range("F48").Value = 177
